$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Two new weekly price records were added to this "Albahaca" dataset.
# They land in the middle of the existing (date-descending-ish) block
# of rows, so the rows below each insertion point shift down by one.
# ------------------------------------------------------------------

# --- Insert #1: new row lands at row 211 --------------------------
$ws.Rows(211).Insert()

$ws.Cells.Item(211, 1).Value = 10
$ws.Cells.Item(211, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(211, 3).Value = "La Araucanía"
$ws.Cells.Item(211, 4).Value = 44748
$ws.Cells.Item(211, 5).Value = 9
$ws.Cells.Item(211, 6).Value = 100112052
$ws.Cells.Item(211, 7).Value = "Albahaca"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 35
$ws.Cells.Item(211, 11).Value = 6000
$ws.Cells.Item(211, 12).Value = 6000
$ws.Cells.Item(211, 13).Value = 6000
$ws.Cells.Item(211, 14).Value = "$/paquete"
$ws.Cells.Item(211, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(211, 16).Value = 6000
$ws.Cells.Item(211, 17).Value = 1
$ws.Cells.Item(211, 18).Value = "Hortaliza"

# --- Insert #2: new row lands at row 231 (after the first shift) --
$ws.Rows(231).Insert()

$ws.Cells.Item(231, 1).Value = 10
$ws.Cells.Item(231, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(231, 3).Value = "La Araucanía"
$ws.Cells.Item(231, 4).Value = 44747
$ws.Cells.Item(231, 5).Value = 9
$ws.Cells.Item(231, 6).Value = 100112052
$ws.Cells.Item(231, 7).Value = "Albahaca"
$ws.Cells.Item(231, 8).Value = "Sin especificar"
$ws.Cells.Item(231, 9).Value = "Primera"
$ws.Cells.Item(231, 10).Value = 50
$ws.Cells.Item(231, 11).Value = 6000
$ws.Cells.Item(231, 12).Value = 6000
$ws.Cells.Item(231, 13).Value = 6000
$ws.Cells.Item(231, 14).Value = "$/paquete"
$ws.Cells.Item(231, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(231, 16).Value = 6000
$ws.Cells.Item(231, 17).Value = 1
$ws.Cells.Item(231, 18).Value = "Hortaliza"
